$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2087.5
$ws.Range("I18").Value = 1415.4
$ws.Range("J18").Value = 2759.6
$ws.Range("K18").Value = 1415.4
$ws.Range("L18").Value = 2759.6
$ws.Range("M18").Value = -1131.4
$ws.Range("N18").Value = -3327.6
$ws.Range("H33").Value = 349.5625
$ws.Range("I33").Value = 342.16666
$ws.Range("K33").Value = 342.16666
$ws.Range("M33").Value = -113.16666
$ws.Range("H39").Value = 89.8
$ws.Range("I39").Value = 99.75
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 299.25
$ws.Range("L39").Value = 150
$ws.Range("M39").Value = -3.25
$ws.Range("N39").Value = -742
$ws.Range("H40").Value = 3843.3333
$ws.Range("J40").Value = 3998.5
$ws.Range("L40").Value = 3998.5
$ws.Range("N40").Value = -4348.5
$ws.Range("H64").Value = 6424.375
$ws.Range("I64").Value = 6066.5
$ws.Range("J64").Value = 7498
$ws.Range("K64").Value = 6066.5
$ws.Range("L64").Value = 7498
$ws.Range("M64").Value = -5818.5
$ws.Range("N64").Value = -7994
$ws.Range("H67").Value = 6424.375
$ws.Range("I67").Value = 6066.5
$ws.Range("J67").Value = 7498
$ws.Range("K67").Value = 6066.5
$ws.Range("L67").Value = 7498
$ws.Range("M67").Value = -5208.5
$ws.Range("N67").Value = -9214
$ws.Range("H100").Value = 2258.1667
$ws.Range("I100").Value = 2077.7778
$ws.Range("K100").Value = 2077.7778
$ws.Range("M100").Value = -1536.7778
$ws.Range("H134").Value = 179999
$ws.Range("J134").Value = 179999
$ws.Range("L134").Value = 179999
$ws.Range("N134").Value = -190139
$ws.Range("H135").Value = 2899
$ws.Range("I135").Value = 799
$ws.Range("K135").Value = 7191
$ws.Range("M135").Value = -4656
$ws.Range("H136").Value = 179694.75
$ws.Range("J136").Value = 179694.75
$ws.Range("L136").Value = 179694.75
$ws.Range("N136").Value = -189894.75
$ws.Range("H137").Value = 1984.75
$ws.Range("I137").Value = 1899
$ws.Range("K137").Value = 5697
$ws.Range("M137").Value = -3147
$ws.Range("H139").Value = 113999
$ws.Range("J139").Value = 113999
$ws.Range("L139").Value = 113999
$ws.Range("N139").Value = -124279
$ws.Range("H140").Value = 99999
$ws.Range("J140").Value = 99999
$ws.Range("L140").Value = 99999
$ws.Range("N140").Value = -110359
$ws.Range("H141").Value = 7570.5
$ws.Range("J141").Value = 12023
$ws.Range("L141").Value = 36069
$ws.Range("N141").Value = -46429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2117.4
$ws.Range("I2").Value = 1997.5
$ws.Range("J2").Value = 2197.3333
$ws.Range("K2").Value = 1997.5
$ws.Range("L2").Value = 2197.3333
$ws.Range("M2").Value = -1884.5
$ws.Range("N2").Value = -2423.3333
$ws.Range("H116").Value = 2117.4
$ws.Range("I116").Value = 1997.5
$ws.Range("J116").Value = 2197.3333
$ws.Range("K116").Value = 1997.5
$ws.Range("L116").Value = 2197.3333
$ws.Range("M116").Value = 296.5
$ws.Range("N116").Value = -6785.3333
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H122").Value = 1205.25
$ws.Range("I122").Value = 1133
$ws.Range("K122").Value = 3399
$ws.Range("M122").Value = -949
$ws.Range("H130").Value = 19849.5
$ws.Range("I130").Value = 19799.666
$ws.Range("K130").Value = 19799.666
$ws.Range("M130").Value = -14779.666
$ws.Range("H132").Value = 2522.4614
$ws.Range("I132").Value = 1466.3334
$ws.Range("J132").Value = 3427.7144
$ws.Range("K132").Value = 4399.0002
$ws.Range("L132").Value = 10283.1432
$ws.Range("M132").Value = -1869.0002
$ws.Range("N132").Value = -15343.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2117.4
$ws.Range("I3").Value = 1997.5
$ws.Range("J3").Value = 2197.3333
$ws.Range("K3").Value = 1997.5
$ws.Range("L3").Value = 2197.3333
$ws.Range("M3").Value = -1883.5
$ws.Range("N3").Value = -2425.3333
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H39").Value = 50000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 50000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 50000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -50778
$ws.Range("H105").Value = 3969.9167
$ws.Range("I105").Value = 2699.8333
$ws.Range("J105").Value = 5240
$ws.Range("K105").Value = 2699.8333
$ws.Range("L105").Value = 5240
$ws.Range("M105").Value = -952.8332999999998
$ws.Range("N105").Value = -8734
$ws.Range("H107").Value = 3464.2
$ws.Range("I107").Value = 3464.2
$ws.Range("K107").Value = 3464.2
$ws.Range("M107").Value = -1544.2
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H131").Value = 69284
$ws.Range("I131").Value = 55000
$ws.Range("J131").Value = 79997
$ws.Range("K131").Value = 55000
$ws.Range("L131").Value = 79997
$ws.Range("M131").Value = -49960
$ws.Range("N131").Value = -90077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 83334910
$ws.Range("I16").Value = 83334910
$ws.Range("K16").Value = 83334910
$ws.Range("M16").Value = -83334623
$ws.Range("H86").Value = 4974.5
$ws.Range("I86").Value = 4828
$ws.Range("K86").Value = 4828
$ws.Range("M86").Value = -3705
$ws.Range("H89").Value = 4974.5
$ws.Range("I89").Value = 4828
$ws.Range("K89").Value = 24140
$ws.Range("M89").Value = -18524
$ws.Range("H113").Value = 83334910
$ws.Range("I113").Value = 83334910
$ws.Range("K113").Value = 83334910
$ws.Range("M113").Value = -83332740
$ws.Range("H134").Value = 3645.2856
$ws.Range("I134").Value = 3671.8462
$ws.Range("K134").Value = 11015.5386
$ws.Range("M134").Value = -8480.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1274.275
$ws.Range("I26").Value = 1402.1143
$ws.Range("J26").Value = 379.4
$ws.Range("K26").Value = 4206.3429
$ws.Range("L26").Value = 1138.2
$ws.Range("M26").Value = -3918.3429
$ws.Range("N26").Value = -1714.2
$ws.Range("H131").Value = 2081.5625
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 238.18182
$ws.Range("I2").Value = 253.85715
$ws.Range("K2").Value = 253.85715
$ws.Range("M2").Value = -140.85715
$ws.Range("H70").Value = 6640.25
$ws.Range("I70").Value = 5504
$ws.Range("J70").Value = 7776.5
$ws.Range("K70").Value = 5504
$ws.Range("L70").Value = 7776.5
$ws.Range("M70").Value = -5234
$ws.Range("N70").Value = -8316.5
$ws.Range("H73").Value = 6640.25
$ws.Range("I73").Value = 5504
$ws.Range("J73").Value = 7776.5
$ws.Range("K73").Value = 5504
$ws.Range("L73").Value = 7776.5
$ws.Range("M73").Value = -4568
$ws.Range("N73").Value = -9648.5
$ws.Range("H80").Value = 1503
$ws.Range("J80").Value = 1503
$ws.Range("L80").Value = 1503
$ws.Range("N80").Value = -3499
$ws.Range("H83").Value = 1503
$ws.Range("J83").Value = 1503
$ws.Range("L83").Value = 7515
$ws.Range("N83").Value = -17499
$ws.Range("H102").Value = 2859.3333
$ws.Range("I102").Value = 3466.5557
$ws.Range("J102").Value = 1037.6666
$ws.Range("K102").Value = 3466.5557
$ws.Range("L102").Value = 1037.6666
$ws.Range("M102").Value = -1844.5557
$ws.Range("N102").Value = -4281.6666
$ws.Range("H113").Value = 2535.5
$ws.Range("I113").Value = 2222.3333
$ws.Range("K113").Value = 2222.3333
$ws.Range("M113").Value = -52.33329999999978
$ws.Range("H122").Value = 2719.3845
$ws.Range("I122").Value = 1843.1428
$ws.Range("J122").Value = 3741.6667
$ws.Range("K122").Value = 5529.428400000001
$ws.Range("L122").Value = 11225.0001
$ws.Range("M122").Value = -3079.428400000001
$ws.Range("N122").Value = -16125.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 860.4
$ws.Range("I93").Value = 893.875
$ws.Range("J93").Value = 726.5
$ws.Range("K93").Value = 893.875
$ws.Range("L93").Value = 726.5
$ws.Range("M93").Value = 354.125
$ws.Range("N93").Value = -3222.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3534449
$ws.Range("I14").Value = 4679265.5
$ws.Range("K14").Value = 4679265.5
$ws.Range("M14").Value = -4679097.5
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H101").Value = 23842
$ws.Range("J101").Value = 23842
$ws.Range("L101").Value = 23842
$ws.Range("N101").Value = -30332
